# Insert two new weekly price rows for "Apio" (Vega Monumental Concepción)
# right before the current row 478, shifting all subsequent data rows down
# by two rows (old row 478 -> new row 480, ..., old row 573 -> new row 575).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 478 (pushes existing 478.. down to 480..)
$ws.Rows("478:479").Insert()

# --- New row 478 (Calidad: Primera) ---
$ws.Cells.Item(478, 1).Value2 = 11
$ws.Cells.Item(478, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(478, 3).Value2 = "Bíobío"
$ws.Cells.Item(478, 4).Value2 = 45244
$ws.Cells.Item(478, 5).Value2 = 8
$ws.Cells.Item(478, 6).Value2 = 100112017
$ws.Cells.Item(478, 7).Value2 = "Apio"
$ws.Cells.Item(478, 8).Value2 = "Americana (o)"
$ws.Cells.Item(478, 9).Value2 = "Primera"
$ws.Cells.Item(478, 10).Value2 = 250
$ws.Cells.Item(478, 11).Value2 = 8000
$ws.Cells.Item(478, 12).Value2 = 9000
$ws.Cells.Item(478, 13).Value2 = 8520
$ws.Cells.Item(478, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(478, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(478, 16).Value2 = 1420
$ws.Cells.Item(478, 17).Value2 = 6
$ws.Cells.Item(478, 18).Value2 = "Hortaliza"

# --- New row 479 (Calidad: Segunda) ---
$ws.Cells.Item(479, 1).Value2 = 11
$ws.Cells.Item(479, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(479, 3).Value2 = "Bíobío"
$ws.Cells.Item(479, 4).Value2 = 45244
$ws.Cells.Item(479, 5).Value2 = 8
$ws.Cells.Item(479, 6).Value2 = 100112017
$ws.Cells.Item(479, 7).Value2 = "Apio"
$ws.Cells.Item(479, 8).Value2 = "Americana (o)"
$ws.Cells.Item(479, 9).Value2 = "Segunda"
$ws.Cells.Item(479, 10).Value2 = 200
$ws.Cells.Item(479, 11).Value2 = 7000
$ws.Cells.Item(479, 12).Value2 = 7000
$ws.Cells.Item(479, 13).Value2 = 7000
$ws.Cells.Item(479, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(479, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(479, 16).Value2 = 1167
$ws.Cells.Item(479, 17).Value2 = 6
$ws.Cells.Item(479, 18).Value2 = "Hortaliza"

# Apply the same date number-format used by column D elsewhere in the sheet
$ws.Cells.Item(478, 4).NumberFormat = $ws.Cells.Item(480, 4).NumberFormat
$ws.Cells.Item(479, 4).NumberFormat = $ws.Cells.Item(481, 4).NumberFormat

Write-Host "Inserted new rows 478-479; dimension now $($ws.UsedRange.Address())"
